# Apply updated odds/values to the weekly FlashScore sheet (Sheet1)
# per the source diff - only the specific cells listed below change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.05  # G2
$ws.Cells.Item(2, 14).Value = 17  # N2
$ws.Cells.Item(2, 19).Value = 2.04  # S2
$ws.Cells.Item(2, 20).Value = 1.86  # T2
# Row 3
$ws.Cells.Item(3, 15).Value = 1.25  # O3
$ws.Cells.Item(3, 16).Value = 4  # P3
$ws.Cells.Item(3, 17).Value = 1.86  # Q3
$ws.Cells.Item(3, 18).Value = 2.04  # R3
$ws.Cells.Item(3, 21).Value = 3  # U3
$ws.Cells.Item(3, 22).Value = 1.4  # V3
# Row 5
$ws.Cells.Item(5, 14).Value = 26  # N5
$ws.Cells.Item(5, 25).Value = 1.75  # Y5
$ws.Cells.Item(5, 26).Value = 2  # Z5
$ws.Cells.Item(5, 39).Value = 67  # AM5
$ws.Cells.Item(5, 43).Value = 51  # AQ5
# Row 7
$ws.Cells.Item(7, 17).Value = 2.6  # Q7
$ws.Cells.Item(7, 18).Value = 1.48  # R7
$ws.Cells.Item(7, 23).Value = 1.57  # W7
$ws.Cells.Item(7, 44).Value = 1.98  # AR7
$ws.Cells.Item(7, 45).Value = 1.88  # AS7
# Row 8
$ws.Cells.Item(8, 7).Value = 2.75  # G8
$ws.Cells.Item(8, 23).Value = 1.44  # W8
$ws.Cells.Item(8, 24).Value = 2.63  # X8
$ws.Cells.Item(8, 29).Value = 11  # AC8
$ws.Cells.Item(8, 30).Value = 29  # AD8
$ws.Cells.Item(8, 31).Value = 23  # AE8
$ws.Cells.Item(8, 40).Value = 10  # AN8
$ws.Cells.Item(8, 42).Value = 21  # AP8
# Row 9
$ws.Cells.Item(9, 17).Value = 1.88  # Q9
$ws.Cells.Item(9, 18).Value = 1.98  # R9
$ws.Cells.Item(9, 21).Value = 3.2  # U9
$ws.Cells.Item(9, 22).Value = 1.36  # V9
$ws.Cells.Item(9, 23).Value = 1.4  # W9
# Row 10
$ws.Cells.Item(10, 7).Value = 1.57  # G10
$ws.Cells.Item(10, 9).Value = 5.75  # I10
$ws.Cells.Item(10, 10).Value = 2.2  # J10
$ws.Cells.Item(10, 23).Value = 1.44  # W10
$ws.Cells.Item(10, 24).Value = 2.63  # X10
$ws.Cells.Item(10, 25).Value = 2.1  # Y10
$ws.Cells.Item(10, 26).Value = 1.67  # Z10
$ws.Cells.Item(10, 28).Value = 6.5  # AB10
$ws.Cells.Item(10, 30).Value = 11  # AD10
$ws.Cells.Item(10, 34).Value = 7.5  # AH10
$ws.Cells.Item(10, 38).Value = 13  # AL10
# Row 11
$ws.Cells.Item(11, 15).Value = 1.3  # O11
$ws.Cells.Item(11, 16).Value = 3.4  # P11
$ws.Cells.Item(11, 17).Value = 2  # Q11
$ws.Cells.Item(11, 18).Value = 1.85  # R11
$ws.Cells.Item(11, 23).Value = 1.4  # W11
# Row 12
$ws.Cells.Item(12, 7).Value = 3  # G12
$ws.Cells.Item(12, 9).Value = 2.45  # I12
$ws.Cells.Item(12, 10).Value = 4  # J12
$ws.Cells.Item(12, 23).Value = 1.57  # W12
$ws.Cells.Item(12, 40).Value = 10  # AN12
$ws.Cells.Item(12, 41).Value = 23  # AO12
# Row 13
$ws.Cells.Item(13, 23).Value = 1.3  # W13
# Row 14
$ws.Cells.Item(14, 7).Value = 1.42  # G14
$ws.Cells.Item(14, 10).Value = 1.91  # J14
$ws.Cells.Item(14, 13).Value = 1.03  # M14
$ws.Cells.Item(14, 15).Value = 1.18  # O14
$ws.Cells.Item(14, 22).Value = 1.5  # V14
$ws.Cells.Item(14, 25).Value = 1.8  # Y14
$ws.Cells.Item(14, 26).Value = 1.95  # Z14
$ws.Cells.Item(14, 28).Value = 7.5  # AB14
# Row 15
$ws.Cells.Item(15, 9).Value = 1.62  # I15
$ws.Cells.Item(15, 11).Value = 2.38  # K15
$ws.Cells.Item(15, 13).Value = 1.01  # M15
$ws.Cells.Item(15, 14).Value = 15  # N15
$ws.Cells.Item(15, 15).Value = 1.2  # O15
$ws.Cells.Item(15, 17).Value = 1.62  # Q15
$ws.Cells.Item(15, 18).Value = 2.25  # R15
$ws.Cells.Item(15, 19).Value = 2  # S15
$ws.Cells.Item(15, 20).Value = 1.8  # T15
$ws.Cells.Item(15, 21).Value = 2.5  # U15
$ws.Cells.Item(15, 22).Value = 1.5  # V15
$ws.Cells.Item(15, 23).Value = 1.3  # W15
$ws.Cells.Item(15, 24).Value = 3.4  # X15
$ws.Cells.Item(15, 25).Value = 1.67  # Y15
$ws.Cells.Item(15, 26).Value = 2.1  # Z15
$ws.Cells.Item(15, 27).Value = 17  # AA15
$ws.Cells.Item(15, 32).Value = 34  # AF15
$ws.Cells.Item(15, 38).Value = 9  # AL15
$ws.Cells.Item(15, 39).Value = 9  # AM15
$ws.Cells.Item(15, 41).Value = 13  # AO15
$ws.Cells.Item(15, 43).Value = 21  # AQ15
# Row 16
$ws.Cells.Item(16, 7).Value = 1.13  # G16
$ws.Cells.Item(16, 10).Value = 1.4  # J16
$ws.Cells.Item(16, 13).Value = 1.01  # M16
$ws.Cells.Item(16, 14).Value = 17  # N16
$ws.Cells.Item(16, 15).Value = 1.07  # O16
$ws.Cells.Item(16, 21).Value = 1.67  # U16
# Row 17
$ws.Cells.Item(17, 9).Value = 1.4  # I17
$ws.Cells.Item(17, 12).Value = 1.83  # L17
$ws.Cells.Item(17, 13).Value = 1.02  # M17
$ws.Cells.Item(17, 14).Value = 12  # N17
$ws.Cells.Item(17, 15).Value = 1.13  # O17
$ws.Cells.Item(17, 16).Value = 5.5  # P17
$ws.Cells.Item(17, 22).Value = 1.67  # V17
# Row 18
$ws.Cells.Item(18, 13).Value = 1.05  # M18
$ws.Cells.Item(18, 15).Value = 1.25  # O18
$ws.Cells.Item(18, 22).Value = 1.36  # V18
# Row 19
$ws.Cells.Item(19, 7).Value = 1.14  # G19
$ws.Cells.Item(19, 15).Value = 1.1  # O19
$ws.Cells.Item(19, 16).Value = 7  # P19
$ws.Cells.Item(19, 21).Value = 1.8  # U19
$ws.Cells.Item(19, 22).Value = 2  # V19
$ws.Cells.Item(19, 25).Value = 2  # Y19
$ws.Cells.Item(19, 26).Value = 1.75  # Z19
$ws.Cells.Item(19, 27).Value = 11  # AA19
$ws.Cells.Item(19, 28).Value = 7.5  # AB19
$ws.Cells.Item(19, 31).Value = 11  # AE19
$ws.Cells.Item(19, 33).Value = 26  # AG19
$ws.Cells.Item(19, 37).Value = 351  # AK19
$ws.Cells.Item(19, 43).Value = 67  # AQ19
# Row 20
$ws.Cells.Item(20, 7).Value = 3.5  # G20
$ws.Cells.Item(20, 8).Value = 3.3  # H20
$ws.Cells.Item(20, 11).Value = 2.1  # K20
$ws.Cells.Item(20, 13).Value = 1.05  # M20
$ws.Cells.Item(20, 15).Value = 1.29  # O20
$ws.Cells.Item(20, 17).Value = 1.98  # Q20
$ws.Cells.Item(20, 18).Value = 1.88  # R20
$ws.Cells.Item(20, 22).Value = 1.33  # V20
$ws.Cells.Item(20, 23).Value = 1.4  # W20
$ws.Cells.Item(20, 24).Value = 2.75  # X20
$ws.Cells.Item(20, 25).Value = 1.73  # Y20
$ws.Cells.Item(20, 26).Value = 2  # Z20
$ws.Cells.Item(20, 27).Value = 11  # AA20
$ws.Cells.Item(20, 28).Value = 17  # AB20
$ws.Cells.Item(20, 31).Value = 29  # AE20
$ws.Cells.Item(20, 33).Value = 10  # AG20
$ws.Cells.Item(20, 38).Value = 8  # AL20
$ws.Cells.Item(20, 39).Value = 10  # AM20
$ws.Cells.Item(20, 40).Value = 9.5  # AN20
# Row 21
$ws.Cells.Item(21, 13).Value = 1.07  # M21
$ws.Cells.Item(21, 15).Value = 1.33  # O21
$ws.Cells.Item(21, 22).Value = 1.25  # V21
# Row 22
$ws.Cells.Item(22, 7).Value = 1.44  # G22
$ws.Cells.Item(22, 13).Value = 1.05  # M22
$ws.Cells.Item(22, 15).Value = 1.3  # O22
$ws.Cells.Item(22, 22).Value = 1.3  # V22
